# Daily attendance processing - 2025-10-19 12:32:56
# Normalises the "last modified by" history stamps in column G: the
# comma-separated list of editors for each affected session row is
# reversed in place (oldest-first -> newest-first ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,11,12,13,14,15,17,21,29,30,31,32,33,38,39,40,41,42,44,48,
          56,57,58,59,60,65,66,67,68,69,71,75,83,84,85,87,88,89,93,95,96,99,
          109,110,111,113,114,115,119,121,122,125,135,136,137,139,140,141,
          145,147,148,151)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $current = [string]$cell.Value2
    $parts = @($current -split ', ')
    $reversed = @($parts[($parts.Count - 1)..0])
    $cell.Value = ($reversed -join ', ')
}
